$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.470.58'
$ws.Range("E2").Value = '  +1.68%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.035.03'
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.14'
$ws.Range("E5").Value = '  +2.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +1.63%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.34'
$ws.Range("E8").Value = '  +3.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0800'
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("E11").Value = '  -1.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.335.98'
$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.44'
$ws.Range("E13").Value = '  +1.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.39'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.745'
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("E16").Value = '  +2.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.040.82'
$ws.Range("E17").Value = '  +1.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.396.48'
$ws.Range("E18").Value = '  +1.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.23'
$ws.Range("E19").Value = '  +0.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.09'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.58'
$ws.Range("E22").Value = '  -0.73%  '

$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("E24").Value = '  +2.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.94'
$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.17'
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("E28").Value = '  +7.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.77'
$ws.Range("E29").Value = '  +0.85%  '

$ws.Range("E30").Value = '  +0.45%  '

$ws.Range("E31").Value = '  +0.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.49'
$ws.Range("E32").Value = '  +1.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0608'
$ws.Range("E33").Value = '  -1.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.03'
$ws.Range("E34").Value = '  +10.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.49'
$ws.Range("E35").Value = '  +1.75%  '

$ws.Range("E36").Value = '  -0.32%  '

$ws.Range("E37").Value = '  +9.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.22'
$ws.Range("E38").Value = '  +2.77%  '

$ws.Range("E39").Value = '  -0.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.471.21'
$ws.Range("E40").Value = '  -1.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0215'
$ws.Range("E41").Value = '  -1.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0930'
$ws.Range("E42").Value = '  +0.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '94.81'
$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("E44").Value = '  +2.39%  '

$ws.Range("E45").Value = '  +17.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.34'
$ws.Range("E46").Value = '  -3.32%  '

$ws.Range("E47").Value = '  -1.28%  '

$ws.Range("E48").Value = '  +1.00%  '

$ws.Range("E49").Value = '  -2.44%  '

$ws.Range("E50").Value = '  +1.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.222.73'
$ws.Range("E51").Value = '  +0.67%  '
